# adicionando funcao de competencia e ano ao gerar os relatorios
# (bug do exel apresentacao nao ser encontrado/preenchido)
#
# Zera todos os valores numericos da tabela de resultados (colunas B:M,
# linhas 2 a 17), incluindo a linha de TOTAL.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:M17").Value = 0
